$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.181.74"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").Value = "2.366.84"
$ws.Range("E3").Value = "  -4.44%  "
$ws.Range("E4").Value = "  +0.95%  "
$ws.Range("D5").Value = "'498.35"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "'129.97"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").Value = "2.372.62"
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("D10").Value = "'0.0968"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "'0.325"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'4.62"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "2.787.79"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "56.122.05"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("D16").Value = "'21.44"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").Value = "2.394.63"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'10.02"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").Value = "'4.00"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").Value = "'306.62"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("D22").Value = "'6.24"
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "'65.12"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "'0.369"
$ws.Range("E26").Value = "  -5.16%  "
$ws.Range("E27").Value = "  -5.05%  "
$ws.Range("D28").Value = "'7.25"
$ws.Range("E28").Value = "  -6.11%  "
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  -8.63%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = "  -7.62%  "
$ws.Range("D36").Value = "'17.54"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").Value = "'1.17"
$ws.Range("E37").Value = "  -7.70%  "
$ws.Range("D38").Value = "'3.76"
$ws.Range("E38").Value = "  -3.98%  "
$ws.Range("D39").Value = "'36.04"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "'0.792"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("D42").Value = "'130.44"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").Value = "'4.74"
$ws.Range("E44").Value = "  -6.73%  "
$ws.Range("D45").Value = "'0.564"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").Value = "'0.0902"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "'240.16"
$ws.Range("E47").Value = "  -9.13%  "
$ws.Range("D48").Value = "'0.0481"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("D50").Value = "'17.00"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'1.54"
$ws.Range("E51").Value = "  -4.20%  "
